# Generate Report for Handback
# Refresh the handback-status report: update the handoff/handback
# timestamps recorded for the zh-cn and de-de target files, and roll the
# "Latest HO Xliff Generate Date" summary on the Overview sheet forward to
# match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: 155a29fd row (row 2) got a fresh handoff + handback cycle.
$zhcn.Range("H2").Value = "2016-09-04 02:52:46"
$zhcn.Range("K2").Value = "2016-09-04 02:53:07"

# de-de sheet: 155a29fd row (row 2) got a fresh handback.
$dede.Range("K2").Value = "2016-09-04 02:53:15"

# Overview sheet: roll the summary date forward for the 155a29fd row.
$overview.Range("G2").Value = "2016-09-04 02:52:50"
